$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 100, pushing existing rows 100-112 down to 101-113
$ws.Rows.Item(100).Insert()

# Populate the newly inserted row 100 with the new record's data
$ws.Range("A100").Value = 5
$ws.Range("B100").Value = "Macroferia Regional de Talca"
$ws.Range("C100").Value = "Maule"
$ws.Range("D100").Value = 45154
$ws.Range("E100").Value = 7
$ws.Range("F100").Value = 100112040
$ws.Range("G100").Value = "Cilantro"
$ws.Range("H100").Value = "Sin especificar"
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 150
$ws.Range("K100").Value = 8000
$ws.Range("L100").Value = 8000
$ws.Range("M100").Value = 8000
$ws.Range("N100").Value = "$/caja 36 atados"
$ws.Range("O100").Value = "Región Metropolitana"
$ws.Range("P100").Value = 222
$ws.Range("Q100").Value = 36
$ws.Range("R100").Value = "Hortaliza"
